# Apply strikethrough formatting to the "Novos resultados ..." paragraph
# (the whole paragraph, including the paragraph mark).
$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Novos resultados")) {
        $target = $p
        break
    }
}

$target.Range.Font.StrikeThrough = 1

# Re-type the " quantidade de deslocamento (plasticidade baixa " span so the
# engine coalesces it into a single run (matching how Word merges runs that
# end up with identical run formatting after an in-place edit).
$rng = $d.Content
$found = $rng.Find.Execute("quantidade de deslocamento (plasticidade baixa ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Select()
    $word.Selection.TypeText("quantidade de deslocamento (plasticidade baixa ")
}

# Move the "_GoBack" bookmark: it used to sit right after
# "Fazer seleção de modelos - bruno"; it now sits in the middle of
# "(Pavel introdução movimento)", between "Pavel i" and "ntrodução".
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Pavel i", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $splitPoint = $rng2.End
    $insertRng = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $insertRng)
}
